$d = $word.ActiveDocument

# Move to the very end of the document content
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd

# Insert a new paragraph after the last one
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)

# Add the hyperlink text + hyperlink object to the new paragraph
$d.Hyperlinks.Add($endRange, "https://mona-media.com/cach-mang-cong-nghiep-4-0-la-gi-loi-ich-hau-qua-va-giai-phap/", $null, $null, "https://mona-media.com/cach-mang-cong-nghiep-4-0-la-gi-loi-ich-hau-qua-va-giai-phap/")

$d.Save()
